# Update LR-pair worksheet with refreshed TPM-derived values.
# Sending/target cluster labels now include the new "ECs" cluster, and the
# per-pair specificity/weight metrics are recomputed against the new TPM data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: ECs -> Fzd10 on FAPs -----------------------------------------
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 0.131499
$ws.Range("H2").Value = 0.394497
$ws.Range("I2").Value = 0.3654391092296077
$ws.Range("J2").Value = 0.3654391092296077
$ws.Range("M2").Value = 0.06694666666666667
$ws.Range("O2").Value = 0.5098924310779488
$ws.Range("P2").Value = 0.5098924310779488
$ws.Range("Q2").Value = 0.00880341972
$ws.Range("R2").Value = 0.07923077748
$ws.Range("S2").Value = 0.1863346358160448
$ws.Range("T2").Value = 0.1863346358160448

# --- Row 3: ECs -> Fzd10 on MuSCs -----------------------------------------
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("G3").Value = 0.131499
$ws.Range("H3").Value = 0.394497
$ws.Range("I3").Value = 0.3654391092296077
$ws.Range("J3").Value = 0.3654391092296077
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.064349
$ws.Range("N3").Value = 0.193047
$ws.Range("O3").Value = 0.4901075689220513
$ws.Range("P3").Value = 0.4901075689220513
$ws.Range("Q3").Value = 0.008461829151000001
$ws.Range("R3").Value = 0.07615646235899999
$ws.Range("S3").Value = 0.179104473413563
$ws.Range("T3").Value = 0.179104473413563

# --- Row 4: FAPs -> Fzd10 on FAPs -----------------------------------------
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("G4").Value = 0.2283393333333333
$ws.Range("H4").Value = 0.685018
$ws.Range("I4").Value = 0.6345608907703922
$ws.Range("J4").Value = 0.6345608907703922
$ws.Range("M4").Value = 0.06694666666666667
$ws.Range("O4").Value = 0.5098924310779488
$ws.Range("P4").Value = 0.5098924310779488
$ws.Range("Q4").Value = 0.01528655723555556
$ws.Range("R4").Value = 0.13757901512
$ws.Range("S4").Value = 0.323557795261904
$ws.Range("T4").Value = 0.323557795261904

# --- Row 5: FAPs -> Fzd10 on MuSCs -----------------------------------------
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 0.2283393333333333
$ws.Range("H5").Value = 0.685018
$ws.Range("I5").Value = 0.6345608907703922
$ws.Range("J5").Value = 0.6345608907703922
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.064349
$ws.Range("N5").Value = 0.193047
$ws.Range("O5").Value = 0.4901075689220513
$ws.Range("P5").Value = 0.4901075689220513
$ws.Range("Q5").Value = 0.01469340776066667
$ws.Range("R5").Value = 0.132240669846
$ws.Range("S5").Value = 0.3110030955084883
$ws.Range("T5").Value = 0.3110030955084883
